# Automatische test-sync: 2025-08-04 21:00:50
# Appends a new test-mail log entry ("Testmail #16") to the "Logs" sheet,
# extends the conditional formatting ranges to cover the new row, and
# bumps the matching category counter on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Determine the last used row before we add the new one.
$oldLastRow = $ws.UsedRange.Rows.Count
$newLastRow = $oldLastRow + 1

# New log row values.
$newRow = @(
    "Wil je dit even doorsturen?",
    "mailmind.test@zohomail.eu",
    "Testmail #16: Wil je dit even doorsturen?",
    "Planning / Afspraak",
    "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl.",
    "2025-08-04 21:00:14",
    "Ja",
    "Ja",
    "Nee",
    "Nee"
)

for ($col = 1; $col -le $newRow.Length; $col++) {
    $ws.Cells.Item($newLastRow, $col).Value = $newRow[$col - 1]
}

# Extend the conditional-formatting sqref ranges (D, G, H, I, J) so they
# keep covering the data from row 2 through the newly added row.
$cfCols = @("D", "G", "H", "I", "J")
foreach ($col in $cfCols) {
    $oldRange = $ws.Range("${col}2:${col}${oldLastRow}")
    $newRange = $ws.Range("${col}2:${col}${newLastRow}")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# Update the Dashboard summary count for "Planning / Afspraak" (7 -> 8).
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Cells.Item(2, 2).Value = 8
